# "CS with scroll and book view"
# The sheet used to hold four rows of ledger data (A1:B4). All but the
# last row are removed, and the remaining row (previously row 4, the
# 2023-07-19 entry) is relocated down to row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the last row's data (A4:B4) down to its new home at A6:B6 first,
# preserving the cell content/type instead of retyping it.
$ws.Range("A4:B4").Cut($ws.Range("A6:B6"))

# Remove the now-stale rows 1-3 that preceded it.
$ws.Range("A1:B3").Clear()
